$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "HTML/CSS" language row is being split into two separate rows:
# the existing row keeps the label but renamed to "HTML", and a brand
# new row is appended for "CSS" with its own scores.
#
# Add the new "CSS" row first (row 22) so that the "CSS" shared string
# is created before the "HTML" shared string (matches the order the
# strings appear in the saved workbook).
$ws.Range("A22").Value = "CSS"
$ws.Range("B22").Value = 0.86
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0.65
$ws.Range("G22").Value = 0.5

# Rename the original combined row (row 7) from "HTML/CSS" to just "HTML".
# Its numeric scores (columns B:G) are unchanged.
$ws.Range("A7").Value = "HTML"

# Update the view: scroll so row 11 is at the top, and select I18.
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
[void]$ws.Range("I18").Select()
